$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 81

# Date/Time/Week columns look numeric to Excel's type inference, so a
# leading apostrophe is used (same as typing them in the Excel UI) to keep
# them as literal text, matching the rest of the column.
$ws.Cells.Item($row, 1).Value = "'2023-06-27"
$ws.Cells.Item($row, 2).Value = "'12:01:33"
$ws.Cells.Item($row, 3).Value = "Tuesday"
$ws.Cells.Item($row, 4).Value = "'26"

$ws.Cells.Item($row, 5).Value = 122831
$ws.Cells.Item($row, 6).Value = 134325
$ws.Cells.Item($row, 7).Value = 163501
$ws.Cells.Item($row, 8).Value = 133794
$ws.Cells.Item($row, 9).Value = 177257
$ws.Cells.Item($row, 10).Value = 115012
$ws.Cells.Item($row, 11).Value = 203681
$ws.Cells.Item($row, 12).Value = 226360
$ws.Cells.Item($row, 13).Value = 176231
$ws.Cells.Item($row, 14).Value = 104336
$ws.Cells.Item($row, 15).Value = 39675
$ws.Cells.Item($row, 16).Value = 33764
$ws.Cells.Item($row, 17).Value = 52226
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 35639
$ws.Cells.Item($row, 20).Value = -1
